# Organized results. Removed old files.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clear out the old 2-cell layout before rebuilding ---
$ws.Range("A1:B2").ClearContents()

# --- New header row (machine names) ---
# Shared-string insertion order matters for the saved sst table, so these
# are written in the same order the original author's tool produced:
# laptop, home, then every row label top-to-bottom, then desktop last.
$ws.Range("B1").Value = "Kerk's laptop"
$ws.Range("C1").Value = "Kerk's home"

# --- Row labels (column A) ---
$ws.Range("A2").Value  = "baseline"
$ws.Range("A3").Value  = "ILAsolveLIN"
$ws.Range("A4").Value  = "ILAsolveVFI"
$ws.Range("A5").Value  = "ILAsolveGSSA"
$ws.Range("A6").Value  = "ILAsimLIN"
$ws.Range("A7").Value  = "ILAsimVFI"
$ws.Range("A8").Value  = "ILAsimGSSA"
$ws.Range("A9").Value  = "OLGsolveLIN"
$ws.Range("A10").Value = "OLGsolveVFI"
$ws.Range("A11").Value = "OLGsolveGSSA"
$ws.Range("A12").Value = "OLGsimLIN"
$ws.Range("A13").Value = "OLGsimVFI"
$ws.Range("A14").Value = "OLGsimGSSA"

$ws.Range("D1").Value = "Kerk's desktop"

# --- Data values ---
$ws.Range("B2").Value = 111.436405553016
$ws.Range("C2").Value = 101.447907824556
$ws.Range("D2").Value = 155.007356226095

$ws.Range("B9").Value  = 0.0254437059920746
$ws.Range("B12").Value = 2163.1343990180098

# --- View: zoom + freeze panes + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 130

$ws.Range("B2").Select()
$win.FreezePanes = $true

$ws.Rows.Item(2).Select()

Write-Output "done"
